$wb = $excel.ActiveWorkbook
$wsRights = $wb.Worksheets.Item("User Rights")
$wsAbout = $wb.Worksheets.Item("About")

foreach ($addr in @("E3","E5","E7","E8","E16","E18","E44","E103","E108","E109","E121")) {
    $wsRights.Range($addr).Value = 'CASE_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E4")) {
    $wsRights.Range($addr).Value = 'CASE_VIEW, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E6")) {
    $wsRights.Range($addr).Value = 'SAMPLE_VIEW, ADDITIONAL_TEST_VIEW, DOCUMENT_VIEW, VISIT_DELETE, CLINICAL_COURSE_VIEW, IMMUNIZATION_VIEW, DOCUMENT_DELETE, IMMUNIZATION_DELETE, TASK_DELETE, TASK_VIEW, PERSON_VIEW, PERSON_DELETE, PATHOGEN_TEST_DELETE, PRESCRIPTION_DELETE, CLINICAL_VISIT_DELETE, CASE_VIEW, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, THERAPY_VIEW, TREATMENT_DELETE'
}

foreach ($addr in @("E9","E10","E11","E12","E13","E14","E15","E17")) {
    $wsRights.Range($addr).Value = 'CASE_VIEW, CASE_EDIT, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E20","E22")) {
    $wsRights.Range($addr).Value = 'IMMUNIZATION_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E21")) {
    $wsRights.Range($addr).Value = 'IMMUNIZATION_VIEW, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E23")) {
    $wsRights.Range($addr).Value = 'VISIT_DELETE, IMMUNIZATION_VIEW, PERSON_VIEW, PERSON_DELETE'
}

foreach ($addr in @("E26")) {
    $wsRights.Range($addr).Value = 'VISIT_DELETE, PERSON_VIEW'
}

foreach ($addr in @("E33")) {
    $wsRights.Range($addr).Value = 'SAMPLE_VIEW, ADDITIONAL_TEST_VIEW, ADDITIONAL_TEST_DELETE, PATHOGEN_TEST_DELETE'
}

foreach ($addr in @("E45","E47","E49","E50","E104")) {
    $wsRights.Range($addr).Value = 'CONTACT_VIEW, CASE_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E46")) {
    $wsRights.Range($addr).Value = 'CONTACT_VIEW, CASE_VIEW, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E48")) {
    $wsRights.Range($addr).Value = 'SAMPLE_VIEW, ADDITIONAL_TEST_VIEW, DOCUMENT_VIEW, VISIT_DELETE, DOCUMENT_DELETE, CONTACT_VIEW, TASK_DELETE, TASK_VIEW, PERSON_VIEW, PERSON_DELETE, PATHOGEN_TEST_DELETE, CASE_VIEW, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE'
}

foreach ($addr in @("E51")) {
    $wsRights.Range($addr).Value = 'CONTACT_VIEW, CASE_VIEW, CONTACT_EDIT, CASE_CREATE, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E52","E53","E54")) {
    $wsRights.Range($addr).Value = 'CONTACT_VIEW, CASE_VIEW, CONTACT_EDIT, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E64")) {
    $wsRights.Range($addr).Value = 'TASK_EDIT, TASK_VIEW'
}

foreach ($addr in @("E67")) {
    $wsRights.Range($addr).Value = 'EVENT_VIEW, DOCUMENT_VIEW, DOCUMENT_DELETE'
}

foreach ($addr in @("E73")) {
    $wsRights.Range($addr).Value = 'EVENTPARTICIPANT_DELETE, SAMPLE_VIEW, ADDITIONAL_TEST_VIEW, EVENT_VIEW, DOCUMENT_VIEW, VISIT_DELETE, DOCUMENT_DELETE, TASK_DELETE, TASK_VIEW, PERSON_VIEW, PERSON_DELETE, PATHOGEN_TEST_DELETE, ACTION_DELETE, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, EVENTPARTICIPANT_VIEW'
}

foreach ($addr in @("E76","E87")) {
    $wsRights.Range($addr).Value = 'EVENT_VIEW, EVENT_EDIT'
}

foreach ($addr in @("E77")) {
    $wsRights.Range($addr).Value = 'EVENT_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E78","E80","E82")) {
    $wsRights.Range($addr).Value = 'EVENT_VIEW, EVENTPARTICIPANT_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E79")) {
    $wsRights.Range($addr).Value = 'EVENT_VIEW, EVENTPARTICIPANT_VIEW, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E81")) {
    $wsRights.Range($addr).Value = 'SAMPLE_VIEW, ADDITIONAL_TEST_VIEW, EVENT_VIEW, VISIT_DELETE, ADDITIONAL_TEST_DELETE, SAMPLE_DELETE, EVENTPARTICIPANT_VIEW, PERSON_VIEW, PERSON_DELETE, PATHOGEN_TEST_DELETE'
}

foreach ($addr in @("E83")) {
    $wsRights.Range($addr).Value = 'EVENT_VIEW, EVENTGROUP_LINK, EVENT_EDIT'
}

foreach ($addr in @("E105")) {
    $wsRights.Range($addr).Value = 'DASHBOARD_CONTACT_VIEW, CONTACT_VIEW, CASE_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E110","E111","E112","E113","E114","E115","E116")) {
    $wsRights.Range($addr).Value = 'CASE_VIEW, THERAPY_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E117","E118","E119","E120")) {
    $wsRights.Range($addr).Value = 'CLINICAL_COURSE_VIEW, CASE_VIEW, THERAPY_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E122")) {
    $wsRights.Range($addr).Value = 'CASE_VIEW, PORT_HEALTH_INFO_VIEW, PERSON_VIEW'
}

foreach ($addr in @("E135")) {
    $wsRights.Range($addr).Value = 'CAMPAIGN_VIEW, CAMPAIGN_FORM_DATA_DELETE, CAMPAIGN_FORM_DATA_VIEW'
}

foreach ($addr in @("E142")) {
    $wsRights.Range($addr).Value = 'TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW'
}

foreach ($addr in @("E143","E145")) {
    $wsRights.Range($addr).Value = 'TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW'
}

foreach ($addr in @("E144")) {
    $wsRights.Range($addr).Value = 'TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW, PERSON_EDIT'
}

foreach ($addr in @("E146")) {
    $wsRights.Range($addr).Value = 'TRAVEL_ENTRY_VIEW, DOCUMENT_VIEW, VISIT_DELETE, DOCUMENT_DELETE, TRAVEL_ENTRY_MANAGEMENT_ACCESS, TASK_DELETE, TASK_VIEW, PERSON_VIEW, PERSON_DELETE'
}

foreach ($addr in @("E151")) {
    $wsRights.Range($addr).Value = 'ENVIRONMENT_SAMPLE_DELETE, ENVIRONMENT_SAMPLE_VIEW, ENVIRONMENT_VIEW, ENVIRONMENT_PATHOGEN_TEST_DELETE'
}

foreach ($addr in @("E152")) {
    $wsRights.Range($addr).Value = 'ENVIRONMENT_CREATE, ENVIRONMENT_VIEW'
}

foreach ($addr in @("E157","E158","E163")) {
    $wsRights.Range($addr).Value = 'ENVIRONMENT_SAMPLE_VIEW, ENVIRONMENT_SAMPLE_EDIT'
}

foreach ($addr in @("E159")) {
    $wsRights.Range($addr).Value = 'ENVIRONMENT_SAMPLE_VIEW, ENVIRONMENT_PATHOGEN_TEST_DELETE'
}

foreach ($addr in @("E185")) {
    $wsRights.Range($addr).Value = 'SAMPLE_VIEW, VISIT_DELETE, SAMPLE_EDIT, PATHOGEN_TEST_EDIT, IMMUNIZATION_VIEW, EXTERNAL_MESSAGE_VIEW, IMMUNIZATION_DELETE, PERSON_VIEW, PATHOGEN_TEST_DELETE, PERSON_EDIT, PERSON_DELETE, SAMPLE_CREATE, PATHOGEN_TEST_CREATE, IMMUNIZATION_EDIT, IMMUNIZATION_CREATE'
}

$wsAbout.Range("A2").Value = '1.0.0'
